$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-09-29 / serial 44468) needs to be inserted
# above the current row 34, pushing the existing rows 34-37 down to 35-38.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44468
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100108
$ws.Range("H34").Value = "Tropicales y subtropicales"
$ws.Range("I34").Value = 100108004
$ws.Range("J34").Value = "Papaya"
$ws.Range("K34").Value = "Cultivar IV Región"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 50
$ws.Range("N34").Value = 16000
$ws.Range("O34").Value = 16000
$ws.Range("P34").Value = 16000
$ws.Range("Q34").Value = "$/bandeja 10 kilos"
$ws.Range("R34").Value = "Provincia del Elquí"
$ws.Range("S34").Value = 1600
$ws.Range("T34").Value = 10
